$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2:E18").Font.Name = "Calibri"

for ($r = 2; $r -le 18; $r++) {
    $ws.Cells.Item($r, 6).Value = "id invalido"
}
